$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: update section header (0017/RESIDENCIAL LAS PILETAS -> 0002/GERENCIA ADMINISTRATIVA) ---
$ws.Range('B8').Value2 = 'GERENCIA ADMINISTRATIVA'
# Force "0002" to be stored as text (keeps leading zeros) without losing A8's original
# number-formatted (General) cell style: stamp Text format long enough to write the
# value, then re-stamp A8's original formatting (still style 4, from the untouched B8)
# back on top so the cell's style index is preserved.
$ws.Range('A8').NumberFormat = '@'
$ws.Range('A8').Value2 = '0002'
$ws.Range('B8').Copy()
$ws.Range('A8').PasteSpecial(-4122)

# --- New employee rows 10-17, text-formatted like the existing data rows (style index 1) ---
# Row 10
$ws.Range('A10').NumberFormat = '@'
$ws.Range('A10').Value2 = '003348'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value2 = 'ROBERTO   HERNANDEZ LOPEZ/SIN UNIFORME'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value2 = '250.00 '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  '
$ws.Range('F10').NumberFormat = '@'
$ws.Range('F10').Value2 = '2025-03-03 00:00:00 '
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value2 = '2030-11-20 00:00:00  '
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value2 = 'F.RETIRO'
$ws.Range('I10').NumberFormat = '@'
$ws.Range('I10').Value2 = '******'
$ws.Range('J10').NumberFormat = '@'
$ws.Range('J10').Value2 = '31/03/04'
$ws.Range('K10').NumberFormat = '@'
$ws.Range('K10').Value2 = '02031000-1'
$ws.Range('L10').NumberFormat = '@'
$ws.Range('L10').Value2 = '68871345'
$ws.Range('M10').NumberFormat = '@'
$ws.Range('M10').Value2 = '0003'
$ws.Range('N10').NumberFormat = '@'
$ws.Range('N10').Value2 = 'GUARDAESPALDA'
$ws.Range('O10').NumberFormat = '@'
$ws.Range('O10').Value2 = '2023'
$ws.Range('P10').NumberFormat = '@'
$ws.Range('P10').Value2 = '0000-00-00'
$ws.Range('Q10').NumberFormat = '@'
$ws.Range('Q10').Value2 = '199680452'
$ws.Range('R10').NumberFormat = '@'
$ws.Range('R10').Value2 = '0614 180868 112 3'
$ws.Range('S10').NumberFormat = '@'
$ws.Range('S10').Value2 = '0368-052358-4'
$ws.Range('T10').NumberFormat = '@'
$ws.Range('T10').Value2 = 'MOTIVO'

# Row 11
$ws.Range('A11').NumberFormat = '@'
$ws.Range('A11').Value2 = '002258'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value2 = 'CELIA NOEMI  BARILLAS ZEPEDA/SIN UNIFORME'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value2 = '0.00 '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  '
$ws.Range('F11').NumberFormat = '@'
$ws.Range('F11').Value2 = '0000-00-00 00:00:00 '
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value2 = '2014-08-16 00:00:00  '
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value2 = 'F.RETIRO'
$ws.Range('I11').NumberFormat = '@'
$ws.Range('I11').Value2 = '******'
$ws.Range('J11').NumberFormat = '@'
$ws.Range('J11').Value2 = '00575460-4'
$ws.Range('K11').NumberFormat = '@'
$ws.Range('K11').Value2 = '269697160000'
$ws.Range('L11').NumberFormat = '@'
$ws.Range('L11').Value2 = '0002'
$ws.Range('M11').NumberFormat = '@'
$ws.Range('M11').Value2 = '2023'
$ws.Range('N11').NumberFormat = '@'
$ws.Range('N11').Value2 = '0000-00-00'
$ws.Range('O11').NumberFormat = '@'
$ws.Range('O11').Value2 = '193732558'
$ws.Range('P11').NumberFormat = '@'
$ws.Range('P11').Value2 = '0816-031173-101-3'
$ws.Range('Q11').NumberFormat = '@'
$ws.Range('Q11').Value2 = '0122-038904-1'
$ws.Range('R11').NumberFormat = '@'
$ws.Range('R11').Value2 = 'MOTIVO'

# Row 12
$ws.Range('A12').NumberFormat = '@'
$ws.Range('A12').Value2 = '010593'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value2 = 'ELIZABETH DEL CARMEN  RECINOS HERNANDEZ/SIN UNIFORME'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value2 = '0.00 '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  '
$ws.Range('F12').NumberFormat = '@'
$ws.Range('F12').Value2 = '2006-01-15 00:00:00 '
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value2 = '2006-01-20 00:00:00  '
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H12').Value2 = 'F.RETIRO'
$ws.Range('I12').NumberFormat = '@'
$ws.Range('I12').Value2 = '******'
$ws.Range('J12').NumberFormat = '@'
$ws.Range('J12').Value2 = '07/01/15'
$ws.Range('K12').NumberFormat = '@'
$ws.Range('K12').Value2 = '02250350-3'
$ws.Range('L12').NumberFormat = '@'
$ws.Range('L12').Value2 = '255127080005'
$ws.Range('M12').NumberFormat = '@'
$ws.Range('M12').Value2 = '0001'
$ws.Range('N12').NumberFormat = '@'
$ws.Range('N12').Value2 = 'GERENTE ADMINISTRATIVO'
$ws.Range('O12').NumberFormat = '@'
$ws.Range('O12').Value2 = '2023'
$ws.Range('P12').NumberFormat = '@'
$ws.Range('P12').Value2 = '0000-00-00'
$ws.Range('Q12').NumberFormat = '@'
$ws.Range('Q12').Value2 = '389694555'
$ws.Range('R12').NumberFormat = '@'
$ws.Range('R12').Value2 = '0614-071169-110-1'
$ws.Range('S12').NumberFormat = '@'
$ws.Range('S12').Value2 = '0322-045110-3'
$ws.Range('T12').NumberFormat = '@'
$ws.Range('T12').Value2 = 'MOTIVO'

# Row 13
$ws.Range('A13').NumberFormat = '@'
$ws.Range('A13').Value2 = '011988'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value2 = 'LUIS FRANCISCO  RECINOS HENRIQUEZ/SIN UNIFORME'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value2 = '0.00 '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  '
$ws.Range('F13').NumberFormat = '@'
$ws.Range('F13').Value2 = '2001-11-17 00:00:00 '
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value2 = '2001-11-17 00:00:00  '
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value2 = 'F.RETIRO'
$ws.Range('I13').NumberFormat = '@'
$ws.Range('I13').Value2 = '******'
$ws.Range('J13').NumberFormat = '@'
$ws.Range('J13').Value2 = '05010802-8'
$ws.Range('K13').NumberFormat = '@'
$ws.Range('K13').Value2 = '0'
$ws.Range('L13').NumberFormat = '@'
$ws.Range('L13').Value2 = '0000'
$ws.Range('M13').NumberFormat = '@'
$ws.Range('M13').Value2 = 'MANTENIMIENTO'
$ws.Range('N13').NumberFormat = '@'
$ws.Range('N13').Value2 = '2023'
$ws.Range('O13').NumberFormat = '@'
$ws.Range('O13').Value2 = '0000-00-00'
$ws.Range('P13').NumberFormat = '@'
$ws.Range('P13').Value2 = '0302-300694-101-0'
$ws.Range('Q13').NumberFormat = '@'
$ws.Range('Q13').Value2 = '0366-043123-5'
$ws.Range('R13').NumberFormat = '@'
$ws.Range('R13').Value2 = 'MOTIVO'

# Row 14
$ws.Range('A14').NumberFormat = '@'
$ws.Range('A14').Value2 = '013203'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value2 = 'JOSE BENJAMIN  SANTOS HERNANDEZ/SIN UNIFORME'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value2 = '305.00 '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '10'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  '
$ws.Range('F14').NumberFormat = '@'
$ws.Range('F14').Value2 = '2017-10-19 00:00:00 '
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value2 = '2017-10-19 00:00:00  '
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value2 = 'F.RETIRO'
$ws.Range('I14').NumberFormat = '@'
$ws.Range('I14').Value2 = '******'
$ws.Range('J14').NumberFormat = '@'
$ws.Range('J14').Value2 = '17/10/19'
$ws.Range('K14').NumberFormat = '@'
$ws.Range('K14').Value2 = '02264938-3'
$ws.Range('L14').NumberFormat = '@'
$ws.Range('L14').Value2 = '275402640001'
$ws.Range('M14').NumberFormat = '@'
$ws.Range('M14').Value2 = '0002'
$ws.Range('N14').NumberFormat = '@'
$ws.Range('N14').Value2 = 'AGENTE DE SEGURIDAD'
$ws.Range('O14').NumberFormat = '@'
$ws.Range('O14').Value2 = '2023'
$ws.Range('P14').NumberFormat = '@'
$ws.Range('P14').Value2 = '0000-00-00'
$ws.Range('Q14').NumberFormat = '@'
$ws.Range('Q14').Value2 = '106750779'
$ws.Range('R14').NumberFormat = '@'
$ws.Range('R14').Value2 = '0210-280575-110-4'
$ws.Range('S14').NumberFormat = '@'
$ws.Range('S14').Value2 = '0311-109439-1'
$ws.Range('T14').NumberFormat = '@'
$ws.Range('T14').Value2 = 'MOTIVO'

# Row 15
$ws.Range('A15').NumberFormat = '@'
$ws.Range('A15').Value2 = '013882'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value2 = 'HUGO LUIS  BARILLAS SANABRIA/SIN UNIFORME'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value2 = '305.00 '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  '
$ws.Range('F15').NumberFormat = '@'
$ws.Range('F15').Value2 = '2021-01-21 00:00:00 '
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value2 = '2021-01-21 00:00:00  '
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value2 = 'F.RETIRO'
$ws.Range('I15').NumberFormat = '@'
$ws.Range('I15').Value2 = '******'
$ws.Range('J15').NumberFormat = '@'
$ws.Range('J15').Value2 = '26/01/21'
$ws.Range('K15').NumberFormat = '@'
$ws.Range('K15').Value2 = '03356690-4'
$ws.Range('L15').NumberFormat = '@'
$ws.Range('L15').Value2 = '313512220018'
$ws.Range('M15').NumberFormat = '@'
$ws.Range('M15').Value2 = '0002'
$ws.Range('N15').NumberFormat = '@'
$ws.Range('N15').Value2 = 'AGENTE DE SEGURIDAD'
$ws.Range('O15').NumberFormat = '@'
$ws.Range('O15').Value2 = '2023'
$ws.Range('P15').NumberFormat = '@'
$ws.Range('P15').Value2 = '0000-00-00'
$ws.Range('Q15').NumberFormat = '@'
$ws.Range('Q15').Value2 = '106856963'
$ws.Range('R15').NumberFormat = '@'
$ws.Range('R15').Value2 = '0210-021185-106-8'
$ws.Range('S15').NumberFormat = '@'
$ws.Range('S15').Value2 = '0322-059358-2'
$ws.Range('T15').NumberFormat = '@'
$ws.Range('T15').Value2 = 'MOTIVO'

# Row 16
$ws.Range('A16').NumberFormat = '@'
$ws.Range('A16').Value2 = '013930'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value2 = 'JUAN ANTONIO  HENRIQUEZ ANZORA/SIN UNIFORME'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value2 = '0.00 '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  '
$ws.Range('F16').NumberFormat = '@'
$ws.Range('F16').Value2 = '2010-02-21 00:00:00 '
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value2 = '2010-02-21 00:00:00  '
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value2 = 'F.RETIRO'
$ws.Range('I16').NumberFormat = '@'
$ws.Range('I16').Value2 = '******'
$ws.Range('J16').NumberFormat = '@'
$ws.Range('J16').Value2 = '19/02/21'
$ws.Range('K16').NumberFormat = '@'
$ws.Range('K16').Value2 = '03536638-4'
$ws.Range('L16').NumberFormat = '@'
$ws.Range('L16').Value2 = '315301870000'
$ws.Range('M16').NumberFormat = '@'
$ws.Range('M16').Value2 = '0001'
$ws.Range('N16').NumberFormat = '@'
$ws.Range('N16').Value2 = 'AGENTE DE SEGURIDAD'
$ws.Range('O16').NumberFormat = '@'
$ws.Range('O16').Value2 = '2023'
$ws.Range('P16').NumberFormat = '@'
$ws.Range('P16').Value2 = '0000-00-00'
$ws.Range('Q16').NumberFormat = '@'
$ws.Range('Q16').Value2 = '107861358'
$ws.Range('R16').NumberFormat = '@'
$ws.Range('R16').Value2 = '0706-300486-101-6'
$ws.Range('S16').NumberFormat = '@'
$ws.Range('S16').Value2 = '0311-255798-3'
$ws.Range('T16').NumberFormat = '@'
$ws.Range('T16').Value2 = 'MOTIVO'

# Row 17
$ws.Range('A17').NumberFormat = '@'
$ws.Range('A17').Value2 = '015517'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value2 = 'JOSE ALBERTO  TOBAR ORTIZ/SIN UNIFORME'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value2 = '0.00 '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  '
$ws.Range('F17').NumberFormat = '@'
$ws.Range('F17').Value2 = '2006-03-23 00:00:00 '
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value2 = '2006-03-23 00:00:00  '
$ws.Range('H17').NumberFormat = '@'
$ws.Range('H17').Value2 = 'F.RETIRO'
$ws.Range('I17').NumberFormat = '@'
$ws.Range('I17').Value2 = '******'
$ws.Range('J17').NumberFormat = '@'
$ws.Range('J17').Value2 = '06/03/23'
$ws.Range('K17').NumberFormat = '@'
$ws.Range('K17').Value2 = '03879639-2'
$ws.Range('L17').NumberFormat = '@'
$ws.Range('L17').Value2 = '321682590002'
$ws.Range('M17').NumberFormat = '@'
$ws.Range('M17').Value2 = '0001'
$ws.Range('N17').NumberFormat = '@'
$ws.Range('N17').Value2 = 'AGENTE DE SEGURIDAD'
$ws.Range('O17').NumberFormat = '@'
$ws.Range('O17').Value2 = '2023'
$ws.Range('P17').NumberFormat = '@'
$ws.Range('P17').Value2 = '0000-00-00'
$ws.Range('Q17').NumberFormat = '@'
$ws.Range('Q17').Value2 = '107885362'
$ws.Range('R17').NumberFormat = '@'
$ws.Range('R17').Value2 = '0309-280188-101-0'
$ws.Range('S17').NumberFormat = '@'
$ws.Range('S17').Value2 = '0322-065675-1'
$ws.Range('T17').NumberFormat = '@'
$ws.Range('T17').Value2 = 'MOTIVO'

